$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Break apart the old merged regions that held the now-removed
#        calendar entries (GS-QC-6301, GS-CC-6208, GS-CC-6202) ---
$ws.Range("C14:C22").UnMerge()
$ws.Range("C14:C20").UnMerge()
$ws.Range("E14:E20").UnMerge()
$ws.Range("B10:B14").UnMerge()
$ws.Range("D10:D14").UnMerge()
$ws.Range("F10:F14").UnMerge()
$ws.Range("B6:B12").UnMerge()
$ws.Range("F6:F12").UnMerge()

# --- 2. Wipe out the old calendar-entry cells (value + formatting) ---
$oldCells = @("B6","F6","B7","F7","B8","F8","B9","F9", `
              "B10","D10","F10","B11","D11","F11","B12","D12","F12", `
              "B13","D13","F13","B14","C14","D14","E14","F14", `
              "C15","E15","C16","E16","C17","E17","C18","E18", `
              "C19","E19","C20","E20","C21")
foreach ($addr in $oldCells) {
    $cell = $ws.Range($addr)
    $cell.Value = ""
    $cell.ClearFormats()
}

# --- 3. Write the new calendar entry: GS-DD-6208, 1:00-2:30, N310 ---
#        It spans Tuesday (C) and Thursday (E), rows 22-28 (1:00-2:30).
$newText = "GS-DD-6208" + [char]160 + "`n1:00-2:30" + [char]160 + "`nN310" + [char]160

$ws.Cells.Item(22,3).Value = $newText
$ws.Cells.Item(22,5).Value = $newText

# --- 4. Apply the existing centered/wrapped style to the new block,
#        matching the original schedule-entry formatting (style index 1) ---
$ws.Cells.Item(22,3).Copy()
$ws.Cells.Item(23,3).PasteSpecial(-4122)
$ws.Cells.Item(24,3).PasteSpecial(-4122)
$ws.Cells.Item(25,3).PasteSpecial(-4122)
$ws.Cells.Item(26,3).PasteSpecial(-4122)
$ws.Cells.Item(27,3).PasteSpecial(-4122)
$ws.Cells.Item(28,3).PasteSpecial(-4122)

$ws.Cells.Item(22,5).Copy()
$ws.Cells.Item(23,5).PasteSpecial(-4122)
$ws.Cells.Item(24,5).PasteSpecial(-4122)
$ws.Cells.Item(25,5).PasteSpecial(-4122)
$ws.Cells.Item(26,5).PasteSpecial(-4122)
$ws.Cells.Item(27,5).PasteSpecial(-4122)
$ws.Cells.Item(28,5).PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- 5. Merge the new block across the 7 rows it occupies ---
$ws.Range("C22:C28").Merge()
$ws.Range("E22:E28").Merge()

# --- 6. Link the location text to a Google Maps driving-directions URL ---
$mapsUrl = "https://www.google.com/maps/dir/Fondren+Gardens,+Houston,+TX/6565+Fannin+St,+Houston,+TX+77030/@29.6290837,-95.5114548,15.31z/data=!4m13!4m12!1m5!1m1!1s0x8640e9093d20ed91:0x82b6198af0aa9bee!2m2!1d-95.5012732!2d29.6263258!1m5!1m1!1s0x8640c071374b0fbd:0x4169184b828fca15!2m2!1d-95.3997225!2d29.7099079"

$ws.Hyperlinks.Add($ws.Cells.Item(28,3), $mapsUrl, "", "", $mapsUrl)
$ws.Hyperlinks.Add($ws.Cells.Item(28,5), $mapsUrl, "", "", $mapsUrl)
